# Applies the commit: split the "事業投資" (business investment) sheet out
# and insert a new "債務" (debt) sheet right before it.
#
# Final sheet order: 土地, 建物, 汽車, 存款, 股票, 其他有價證券, 債務, 事業投資

$wb = $excel.ActiveWorkbook

$lastSheet  = $wb.Worksheets.Item($wb.Worksheets.Count)   # 事業投資 (currently last)
$headerSrc  = $wb.Worksheets.Item(1).Range("B1")            # a cell carrying the bold/bordered header style
$dataSrc    = $wb.Worksheets.Item(1).Range("B2")            # a cell carrying the plain data style

# Insert the new sheet immediately before "事業投資" so the tab order becomes
# ... 其他有價證券, 債務, 事業投資
$debt = $wb.Worksheets.Add($lastSheet)
$debt.Name = "債務"

# ---- header row -----------------------------------------------------
$headers = @{
    "B1" = "species";
    "C1" = "debtor";
    "D1" = "owner";
    "E1" = "total";
    "F1" = "register_date";
    "G1" = "register_reason";
    "H1" = "property_category";
    "I1" = "category";
    "J1" = "date";
    "K1" = "legislator_name";
    "L1" = "legislator_id";
    "M1" = "source_file";
    "N1" = "index";
}
foreach ($addr in $headers.Keys) {
    $debt.Range($addr).Value = $headers[$addr]
}

# ---- data rows --------------------------------------------------------
$rows = @(
    @{
        "A" = 127; "B" = "長期擔保貸款"; "C" = "葉宜津";
        "D" = "台灣土地銀行民權分行臺北市中山區民權西路"; "E" = 1211209;
        "F" = "88年04月07日"; "G" = "88年4月7曰"; "H" = "debt";
        "I" = "normal"; "J" = "2011-11-22"; "K" = "葉宜津";
        "L" = 855; "M" = "tmp14431"; "N" = 127
    },
    @{
        "A" = 128; "B" = "長期擔保貸款"; "C" = "趙哲宏";
        "D" = "台灣土地豳行新營分行臺南市新營區新進路"; "E" = 4000000;
        "F" = "100年04月15日"; "G" = "100年4月15日"; "H" = "debt";
        "I" = "normal"; "J" = "2011-11-22"; "K" = "葉宜津";
        "L" = 855; "M" = "tmp14431"; "N" = 128
    }
)

$rowIndex = 2
foreach ($row in $rows) {
    foreach ($col in $row.Keys) {
        $addr = "$col$rowIndex"
        $val = $row[$col]
        if ($col -eq "J") {
            # "2011-11-22" typed literally gets auto-converted to a date serial
            # by the COM value setter. Route it through a formula + paste-values
            # round-trip so it lands as the plain text shared string instead.
            $debt.Range($addr).Formula = "=""$val"""
            $debt.Range($addr).Copy() | Out-Null
            $debt.Range($addr).PasteSpecial(-4163) | Out-Null
        } else {
            $debt.Range($addr).Value = $val
        }
    }
    $rowIndex++
}

# ---- styling ------------------------------------------------------------
# Header row: bold + thin border + centered (matches the other sheets).
# (Column A has no header label on any sheet in this workbook, so leave A1 alone.)
$headerSrc.Copy() | Out-Null
$debt.Range("B1:N1").PasteSpecial(-4122) | Out-Null

# Column A (index numbers) on the data rows carries the same style as the
# header, matching the other sheets in this workbook.
$headerSrc.Copy() | Out-Null
$debt.Range("A2:A3").PasteSpecial(-4122) | Out-Null

# Remaining data cells: plain style.
$dataSrc.Copy() | Out-Null
$debt.Range("B2:N3").PasteSpecial(-4122) | Out-Null

# The "index" column on the pre-existing "事業投資" sheet shifts from 134 to
# 133 now that the new debt rows (127, 128) have been spliced in ahead of it.
$invest = $wb.Worksheets.Item($wb.Worksheets.Count)
$invest.Range("A2").Value = 133
$invest.Range("N2").Value = 133

$debt.Range("A1").Select()
$excel.CutCopyMode = $false
